$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 17 (pushes the old row 17 data row etc. down by one)
$ws.Rows("17:17").Insert()

# Give the new header-like row 17 the same look as row 16 (bold, centered,
# middle-aligned, wrapped, bordered) by copying the formats only from row 16
$ws.Range("A16:N16").Copy()
$ws.Range("A17:N17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Make the new row a bit shorter than the real header row and fill it with
# the column index numbers 1..14
$ws.Rows("17:17").RowHeight = 16.2
for ($c = 1; $c -le 14; $c++) {
    $ws.Cells.Item(17, $c).Value = $c
}

# Turn the new row into the table header for an AutoFilter
$ws.Range("A17:N17").AutoFilter()

# Register the (hidden, sheet-scoped) _FilterDatabase defined name that
# Excel creates when AutoFilter is applied interactively
$ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$17:`$N`$17")
$fdName = $ws.Names.Item($ws.Names.Count)
$fdName.Visible = $false

# Restore/update the active selection to match the saved view
$ws.Range("H24").Select()
